# Append three new student records to the table on Sheet1 (rows 25-27),
# continuing the existing ID sequence (23 -> 24, 25, 26).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Each inner array holds one row's values in column order:
# ID, Name, Father's Name, DOB, Place of Birth, Profession, Religion, Sex,
# Marital Status, Blood Group, Email, Mobile
$newRows = @(
    @("24", "Md. Masuk Al Hussain ", "Tarik Hossain", "2023-05-28", "Rajshahi", "Software Engineer", "Islam", "Male", "Single", "A+", "masukalhussain3@gmail.com", "01709014797"),
    @("25", "Ziyana Islam", "Quamrul Islam", "2019-01-10", "Dhaka, Bangladesh", "Student", "Islam", "Female", "Single", "A+", "ziyana@gmail.com", "01709014797"),
    @("26", "Affan Islam", "Salman Islam", "2020-01-17", "Dhaka, Bangladesh", "Student", "Islam", "Male", "Single", "A+", "affan@gmail.com", "01709014797")
)

$startRow = 25

# Columns whose values look numeric/date-like and must be forced to text so
# Excel doesn't silently convert them to a number/date (ID column, the
# DOB text-dates, and mobile numbers with a leading zero).
$textColumns = @(1, 4, 12)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $values = $newRows[$i]
    for ($col = 1; $col -le $values.Count; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        if ($textColumns -contains $col) {
            # Mark as Text before assigning so "24", "2023-05-28" and
            # "01709014797" are kept verbatim instead of becoming a
            # number/date serial, then drop back to the default "Normal"
            # style so the cell doesn't carry a one-off text format (matches
            # the plain, unstyled cells used throughout the rest of the
            # sheet).
            $cell.NumberFormat = "@"
            $cell.Value = $values[$col - 1]
            $cell.Style = "Normal"
        } else {
            $cell.Value = $values[$col - 1]
        }
    }
}
